$wb = $excel.ActiveWorkbook

# Rename the "Requested quantity" header cells on the existing sheets
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet as the last tab
$wsForecast = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsForecast.Name = "PO Forecast"

# Reuse the existing bold/centered/bordered header style from "Weekly Quantity" A1:B1
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Reuse the existing date-formatted style from "Weekly Quantity" A2 for column A
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A23").PasteSpecial(-4122)

$headers = @("ds", "PO_Forecast", "yhat_lower", "yhat_upper")
for ($c = 0; $c -lt 4; $c++) {
  $wsForecast.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$forecastRows = @(
  @(45361.99999999999, 7, -58.66251824505196, 75.18724820360687),
  @(45375.99999999999, 14, -52.80064889116798, 82.93226948269242),
  @(45382.99999999999, 17, -51.69030784395048, 80.93990520226912),
  @(45459.99999999999, 51, -15.83151401784165, 120.7927419503395),
  @(45466.99999999999, 54, -14.82851113211432, 121.467085084084),
  @(45480.99999999999, 60, -7.758247576342394, 127.8225090111323),
  @(45487.99999999999, 63, -6.651061366440854, 128.205489567631),
  @(45508.99999999999, 72, 3.654180590770799, 144.9323715427751),
  @(45515.99999999999, 75, 9.970425161821199, 141.031431559562),
  @(45522.99999999999, 79, 12.53342136649112, 145.3696614669753),
  @(45543.99999999999, 88, 21.45974138903028, 153.9962998964121),
  @(45550.99999999999, 91, 24.18799760965796, 157.8810774956781),
  @(45585.99999999999, 106, 34.86093868487427, 173.7354266303738),
  @(45592.99999999999, 109, 39.49214159507422, 177.7672580779509),
  @(45599.99999999999, 113, 44.73067735462391, 176.2443977888817),
  @(45606.99999999999, 116, 47.4283844759669, 182.5575919284828),
  @(45613.99999999999, 119, 56.00704601608413, 190.4562785869499),
  @(45620.99999999999, 122, 47.34439673782369, 186.8279510588012),
  @(45627.99999999999, 125, 63.7279721713165, 190.1377855371659),
  @(45634.99999999999, 128, 56.738949741182, 196.9726630955518),
  @(45641.99999999999, 131, 62.32690460737446, 201.9913091576477),
  @(45648.99999999999, 134, 67.45829303339501, 198.2506350745705)
)

for ($r = 0; $r -lt $forecastRows.Count; $r++) {
  for ($c = 0; $c -lt 4; $c++) {
    $wsForecast.Cells.Item($r + 2, $c + 1).Value = $forecastRows[$r][$c]
  }
}
